$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "F0481-TGAGATCGAA"
$ws.Range("C2").Value = "AATGATACGGCGACCACCGAGATCTACACTGAGATCGAATCGTCGGCAGCGTC"
$ws.Range("B3").Value = "F0482-AGTGATCACC"
$ws.Range("C3").Value = "AATGATACGGCGACCACCGAGATCTACACAGTGATCACCTCGTCGGCAGCGTC"
$ws.Range("B4").Value = "F0483-AGGAACAGAT"
$ws.Range("C4").Value = "AATGATACGGCGACCACCGAGATCTACACAGGAACAGATTCGTCGGCAGCGTC"
$ws.Range("B5").Value = "F0484-AGTGATGGAA"
$ws.Range("C5").Value = "AATGATACGGCGACCACCGAGATCTACACAGTGATGGAATCGTCGGCAGCGTC"
$ws.Range("B6").Value = "F0485-AGATGCTGCT"
$ws.Range("C6").Value = "AATGATACGGCGACCACCGAGATCTACACAGATGCTGCTTCGTCGGCAGCGTC"
$ws.Range("B7").Value = "F0486-ACGTCTGGTA"
$ws.Range("C7").Value = "AATGATACGGCGACCACCGAGATCTACACACGTCTGGTATCGTCGGCAGCGTC"
$ws.Range("B8").Value = "F0487-ACTGGACAAG"
$ws.Range("C8").Value = "AATGATACGGCGACCACCGAGATCTACACACTGGACAAGTCGTCGGCAGCGTC"
$ws.Range("B9").Value = "F0488-GAGAGTCAGT"
$ws.Range("C9").Value = "AATGATACGGCGACCACCGAGATCTACACGAGAGTCAGTTCGTCGGCAGCGTC"
$ws.Range("B10").Value = "F0489-AACACAAGCA"
$ws.Range("C10").Value = "AATGATACGGCGACCACCGAGATCTACACAACACAAGCATCGTCGGCAGCGTC"
$ws.Range("B11").Value = "F0490-CTCAGTCCTA"
$ws.Range("C11").Value = "AATGATACGGCGACCACCGAGATCTACACCTCAGTCCTATCGTCGGCAGCGTC"
$ws.Range("B12").Value = "F0491-AAGATGCAGG"
$ws.Range("C12").Value = "AATGATACGGCGACCACCGAGATCTACACAAGATGCAGGTCGTCGGCAGCGTC"
$ws.Range("B13").Value = "F0492-AACCACACGT"
$ws.Range("C13").Value = "AATGATACGGCGACCACCGAGATCTACACAACCACACGTTCGTCGGCAGCGTC"
$ws.Range("B14").Value = "F0493-TACACCACCT"
$ws.Range("C14").Value = "AATGATACGGCGACCACCGAGATCTACACTACACCACCTTCGTCGGCAGCGTC"
$ws.Range("B15").Value = "F0494-AGCTTGTTCT"
$ws.Range("C15").Value = "AATGATACGGCGACCACCGAGATCTACACAGCTTGTTCTTCGTCGGCAGCGTC"
$ws.Range("B16").Value = "F0495-CCTGAACAGA"
$ws.Range("C16").Value = "AATGATACGGCGACCACCGAGATCTACACCCTGAACAGATCGTCGGCAGCGTC"
$ws.Range("B17").Value = "F0496-GCTACTAGTG"
$ws.Range("C17").Value = "AATGATACGGCGACCACCGAGATCTACACGCTACTAGTGTCGTCGGCAGCGTC"
$ws.Range("B18").Value = "F0497-CGTCGATCAT"
$ws.Range("C18").Value = "AATGATACGGCGACCACCGAGATCTACACCGTCGATCATTCGTCGGCAGCGTC"
$ws.Range("B19").Value = "F0498-GATGTCAGAC"
$ws.Range("C19").Value = "AATGATACGGCGACCACCGAGATCTACACGATGTCAGACTCGTCGGCAGCGTC"
$ws.Range("B20").Value = "F0499-GTGTACATCC"
$ws.Range("C20").Value = "AATGATACGGCGACCACCGAGATCTACACGTGTACATCCTCGTCGGCAGCGTC"
$ws.Range("B21").Value = "F0500-ACATGGTTGG"
$ws.Range("C21").Value = "AATGATACGGCGACCACCGAGATCTACACACATGGTTGGTCGTCGGCAGCGTC"
$ws.Range("B22").Value = "F0501-GGTCAACGAA"
$ws.Range("C22").Value = "AATGATACGGCGACCACCGAGATCTACACGGTCAACGAATCGTCGGCAGCGTC"
$ws.Range("B23").Value = "F0502-AGAGTCTGAT"
$ws.Range("C23").Value = "AATGATACGGCGACCACCGAGATCTACACAGAGTCTGATTCGTCGGCAGCGTC"
$ws.Range("B24").Value = "F0503-TGGTCACTTC"
$ws.Range("C24").Value = "AATGATACGGCGACCACCGAGATCTACACTGGTCACTTCTCGTCGGCAGCGTC"
$ws.Range("B25").Value = "F0504-TCAACTGTCA"
$ws.Range("C25").Value = "AATGATACGGCGACCACCGAGATCTACACTCAACTGTCATCGTCGGCAGCGTC"
$ws.Range("B26").Value = "F0505-AGCACAGGAT"
$ws.Range("C26").Value = "AATGATACGGCGACCACCGAGATCTACACAGCACAGGATTCGTCGGCAGCGTC"
$ws.Range("B27").Value = "F0506-CATCCTACAG"
$ws.Range("C27").Value = "AATGATACGGCGACCACCGAGATCTACACCATCCTACAGTCGTCGGCAGCGTC"
$ws.Range("B28").Value = "F0507-TGTACGTCGA"
$ws.Range("C28").Value = "AATGATACGGCGACCACCGAGATCTACACTGTACGTCGATCGTCGGCAGCGTC"
$ws.Range("B29").Value = "F0508-AACAAGGAAG"
$ws.Range("C29").Value = "AATGATACGGCGACCACCGAGATCTACACAACAAGGAAGTCGTCGGCAGCGTC"
$ws.Range("B30").Value = "F0509-GTAGAACCAG"
$ws.Range("C30").Value = "AATGATACGGCGACCACCGAGATCTACACGTAGAACCAGTCGTCGGCAGCGTC"
$ws.Range("B31").Value = "F0510-CTTGGTAGAG"
$ws.Range("C31").Value = "AATGATACGGCGACCACCGAGATCTACACCTTGGTAGAGTCGTCGGCAGCGTC"
$ws.Range("B32").Value = "F0511-GATCTCCACA"
$ws.Range("C32").Value = "AATGATACGGCGACCACCGAGATCTACACGATCTCCACATCGTCGGCAGCGTC"
$ws.Range("B33").Value = "F0512-GACAAGTCGT"
$ws.Range("C33").Value = "AATGATACGGCGACCACCGAGATCTACACGACAAGTCGTTCGTCGGCAGCGTC"
$ws.Range("B34").Value = "F0513-AGAAGACCTA"
$ws.Range("C34").Value = "AATGATACGGCGACCACCGAGATCTACACAGAAGACCTATCGTCGGCAGCGTC"
$ws.Range("B35").Value = "F0514-CTCTCAGATC"
$ws.Range("C35").Value = "AATGATACGGCGACCACCGAGATCTACACCTCTCAGATCTCGTCGGCAGCGTC"
$ws.Range("B36").Value = "F0515-TTGCACACTC"
$ws.Range("C36").Value = "AATGATACGGCGACCACCGAGATCTACACTTGCACACTCTCGTCGGCAGCGTC"
$ws.Range("B37").Value = "F0516-CATGGAGCTA"
$ws.Range("C37").Value = "AATGATACGGCGACCACCGAGATCTACACCATGGAGCTATCGTCGGCAGCGTC"
$ws.Range("B38").Value = "F0517-GATGTTCCAT"
$ws.Range("C38").Value = "AATGATACGGCGACCACCGAGATCTACACGATGTTCCATTCGTCGGCAGCGTC"
$ws.Range("B39").Value = "F0518-CTTGTAGACG"
$ws.Range("C39").Value = "AATGATACGGCGACCACCGAGATCTACACCTTGTAGACGTCGTCGGCAGCGTC"
$ws.Range("B40").Value = "F0519-ATGTCTAGAC"
$ws.Range("C40").Value = "AATGATACGGCGACCACCGAGATCTACACATGTCTAGACTCGTCGGCAGCGTC"
$ws.Range("B41").Value = "F0520-GAGAGTTGCT"
$ws.Range("C41").Value = "AATGATACGGCGACCACCGAGATCTACACGAGAGTTGCTTCGTCGGCAGCGTC"
$ws.Range("B42").Value = "F0521-TTGAACTAGC"
$ws.Range("C42").Value = "AATGATACGGCGACCACCGAGATCTACACTTGAACTAGCTCGTCGGCAGCGTC"
$ws.Range("B43").Value = "F0522-GCATGCAAGA"
$ws.Range("C43").Value = "AATGATACGGCGACCACCGAGATCTACACGCATGCAAGATCGTCGGCAGCGTC"
$ws.Range("B44").Value = "F0523-CTACAGTACC"
$ws.Range("C44").Value = "AATGATACGGCGACCACCGAGATCTACACCTACAGTACCTCGTCGGCAGCGTC"
$ws.Range("B45").Value = "F0524-TCTAGTGCAG"
$ws.Range("C45").Value = "AATGATACGGCGACCACCGAGATCTACACTCTAGTGCAGTCGTCGGCAGCGTC"
$ws.Range("B46").Value = "F0525-CAACTGAAGG"
$ws.Range("C46").Value = "AATGATACGGCGACCACCGAGATCTACACCAACTGAAGGTCGTCGGCAGCGTC"
$ws.Range("B47").Value = "F0526-CCTTCGTGAT"
$ws.Range("C47").Value = "AATGATACGGCGACCACCGAGATCTACACCCTTCGTGATTCGTCGGCAGCGTC"
$ws.Range("B48").Value = "F0527-TCTAGGTCTT"
$ws.Range("C48").Value = "AATGATACGGCGACCACCGAGATCTACACTCTAGGTCTTTCGTCGGCAGCGTC"
$ws.Range("B49").Value = "F0528-GTCCACTAGA"
$ws.Range("C49").Value = "AATGATACGGCGACCACCGAGATCTACACGTCCACTAGATCGTCGGCAGCGTC"
$ws.Range("B50").Value = "F0529-AAGAAGAGTG"
$ws.Range("C50").Value = "AATGATACGGCGACCACCGAGATCTACACAAGAAGAGTGTCGTCGGCAGCGTC"
$ws.Range("B51").Value = "F0530-GATGCTCATC"
$ws.Range("C51").Value = "AATGATACGGCGACCACCGAGATCTACACGATGCTCATCTCGTCGGCAGCGTC"
$ws.Range("B52").Value = "F0531-TAGAACTGAC"
$ws.Range("C52").Value = "AATGATACGGCGACCACCGAGATCTACACTAGAACTGACTCGTCGGCAGCGTC"
$ws.Range("B53").Value = "F0532-TCCATGTTCG"
$ws.Range("C53").Value = "AATGATACGGCGACCACCGAGATCTACACTCCATGTTCGTCGTCGGCAGCGTC"
$ws.Range("B54").Value = "F0533-GGTCTTGCTT"
$ws.Range("C54").Value = "AATGATACGGCGACCACCGAGATCTACACGGTCTTGCTTTCGTCGGCAGCGTC"
$ws.Range("B55").Value = "F0534-ATCGTTCTGA"
$ws.Range("C55").Value = "AATGATACGGCGACCACCGAGATCTACACATCGTTCTGATCGTCGGCAGCGTC"
$ws.Range("B56").Value = "F0535-GCATCATCTG"
$ws.Range("C56").Value = "AATGATACGGCGACCACCGAGATCTACACGCATCATCTGTCGTCGGCAGCGTC"
$ws.Range("B57").Value = "F0536-CCAACTTGAA"
$ws.Range("C57").Value = "AATGATACGGCGACCACCGAGATCTACACCCAACTTGAATCGTCGGCAGCGTC"
$ws.Range("B58").Value = "F0537-TGTAGACAGT"
$ws.Range("C58").Value = "AATGATACGGCGACCACCGAGATCTACACTGTAGACAGTTCGTCGGCAGCGTC"
$ws.Range("B59").Value = "F0538-CAAGAAGGTT"
$ws.Range("C59").Value = "AATGATACGGCGACCACCGAGATCTACACCAAGAAGGTTTCGTCGGCAGCGTC"
$ws.Range("B60").Value = "F0539-TTGTTCGAGA"
$ws.Range("C60").Value = "AATGATACGGCGACCACCGAGATCTACACTTGTTCGAGATCGTCGGCAGCGTC"
$ws.Range("B61").Value = "F0540-CTCTACATCA"
$ws.Range("C61").Value = "AATGATACGGCGACCACCGAGATCTACACCTCTACATCATCGTCGGCAGCGTC"
$ws.Range("B62").Value = "F0541-TTCCACGTCT"
$ws.Range("C62").Value = "AATGATACGGCGACCACCGAGATCTACACTTCCACGTCTTCGTCGGCAGCGTC"
$ws.Range("B63").Value = "F0542-TACCTTCACA"
$ws.Range("C63").Value = "AATGATACGGCGACCACCGAGATCTACACTACCTTCACATCGTCGGCAGCGTC"
$ws.Range("B64").Value = "F0543-CTACCAAGAG"
$ws.Range("C64").Value = "AATGATACGGCGACCACCGAGATCTACACCTACCAAGAGTCGTCGGCAGCGTC"
$ws.Range("B65").Value = "F0544-CATGGATCAC"
$ws.Range("C65").Value = "AATGATACGGCGACCACCGAGATCTACACCATGGATCACTCGTCGGCAGCGTC"
$ws.Range("B66").Value = "F0545-GAACTGACAA"
$ws.Range("C66").Value = "AATGATACGGCGACCACCGAGATCTACACGAACTGACAATCGTCGGCAGCGTC"
$ws.Range("B67").Value = "F0546-TTGAGACCTT"
$ws.Range("C67").Value = "AATGATACGGCGACCACCGAGATCTACACTTGAGACCTTTCGTCGGCAGCGTC"
$ws.Range("B68").Value = "F0547-TGCACCTTCA"
$ws.Range("C68").Value = "AATGATACGGCGACCACCGAGATCTACACTGCACCTTCATCGTCGGCAGCGTC"
$ws.Range("B69").Value = "F0548-AACGTGACAT"
$ws.Range("C69").Value = "AATGATACGGCGACCACCGAGATCTACACAACGTGACATTCGTCGGCAGCGTC"
$ws.Range("B70").Value = "F0549-CCTGACAGAA"
$ws.Range("C70").Value = "AATGATACGGCGACCACCGAGATCTACACCCTGACAGAATCGTCGGCAGCGTC"
$ws.Range("B71").Value = "F0550-CAAGCTAGCT"
$ws.Range("C71").Value = "AATGATACGGCGACCACCGAGATCTACACCAAGCTAGCTTCGTCGGCAGCGTC"
$ws.Range("B72").Value = "F0551-CTAGTGAGAG"
$ws.Range("C72").Value = "AATGATACGGCGACCACCGAGATCTACACCTAGTGAGAGTCGTCGGCAGCGTC"
$ws.Range("B73").Value = "F0552-TCTCTACCTG"
$ws.Range("C73").Value = "AATGATACGGCGACCACCGAGATCTACACTCTCTACCTGTCGTCGGCAGCGTC"
$ws.Range("B74").Value = "F0553-ACTGTGTTCG"
$ws.Range("C74").Value = "AATGATACGGCGACCACCGAGATCTACACACTGTGTTCGTCGTCGGCAGCGTC"
$ws.Range("B75").Value = "F0554-AAGCTACTCG"
$ws.Range("C75").Value = "AATGATACGGCGACCACCGAGATCTACACAAGCTACTCGTCGTCGGCAGCGTC"
$ws.Range("B76").Value = "F0555-AACTCGATGA"
$ws.Range("C76").Value = "AATGATACGGCGACCACCGAGATCTACACAACTCGATGATCGTCGGCAGCGTC"
$ws.Range("B77").Value = "F0556-TAGCAGACCT"
$ws.Range("C77").Value = "AATGATACGGCGACCACCGAGATCTACACTAGCAGACCTTCGTCGGCAGCGTC"
$ws.Range("B78").Value = "F0557-TGAGGACGTA"
$ws.Range("C78").Value = "AATGATACGGCGACCACCGAGATCTACACTGAGGACGTATCGTCGGCAGCGTC"
$ws.Range("B79").Value = "F0558-AGTCGAACTA"
$ws.Range("C79").Value = "AATGATACGGCGACCACCGAGATCTACACAGTCGAACTATCGTCGGCAGCGTC"
$ws.Range("B80").Value = "F0559-GTAGCTACGT"
$ws.Range("C80").Value = "AATGATACGGCGACCACCGAGATCTACACGTAGCTACGTTCGTCGGCAGCGTC"
$ws.Range("B81").Value = "F0560-AGGATCTGAG"
$ws.Range("C81").Value = "AATGATACGGCGACCACCGAGATCTACACAGGATCTGAGTCGTCGGCAGCGTC"
$ws.Range("B82").Value = "F0561-CTACAAGTAG"
$ws.Range("C82").Value = "AATGATACGGCGACCACCGAGATCTACACCTACAAGTAGTCGTCGGCAGCGTC"
$ws.Range("B83").Value = "F0562-AGATGTGGAA"
$ws.Range("C83").Value = "AATGATACGGCGACCACCGAGATCTACACAGATGTGGAATCGTCGGCAGCGTC"
$ws.Range("B84").Value = "F0563-TACTACAGCT"
$ws.Range("C84").Value = "AATGATACGGCGACCACCGAGATCTACACTACTACAGCTTCGTCGGCAGCGTC"
$ws.Range("B85").Value = "F0564-GCTCATCAAC"
$ws.Range("C85").Value = "AATGATACGGCGACCACCGAGATCTACACGCTCATCAACTCGTCGGCAGCGTC"
$ws.Range("B86").Value = "F0565-TTCAGACGTA"
$ws.Range("C86").Value = "AATGATACGGCGACCACCGAGATCTACACTTCAGACGTATCGTCGGCAGCGTC"
$ws.Range("B87").Value = "F0566-ACAAGGTGCA"
$ws.Range("C87").Value = "AATGATACGGCGACCACCGAGATCTACACACAAGGTGCATCGTCGGCAGCGTC"
$ws.Range("B88").Value = "F0567-TGGAGAGATC"
$ws.Range("C88").Value = "AATGATACGGCGACCACCGAGATCTACACTGGAGAGATCTCGTCGGCAGCGTC"
$ws.Range("B89").Value = "F0568-GGACTCTCTA"
$ws.Range("C89").Value = "AATGATACGGCGACCACCGAGATCTACACGGACTCTCTATCGTCGGCAGCGTC"
$ws.Range("B90").Value = "F0569-GAGCAGTTGA"
$ws.Range("C90").Value = "AATGATACGGCGACCACCGAGATCTACACGAGCAGTTGATCGTCGGCAGCGTC"
$ws.Range("B91").Value = "F0570-ACGTCGAAGT"
$ws.Range("C91").Value = "AATGATACGGCGACCACCGAGATCTACACACGTCGAAGTTCGTCGGCAGCGTC"
$ws.Range("B92").Value = "F0571-TGAACAACAC"
$ws.Range("C92").Value = "AATGATACGGCGACCACCGAGATCTACACTGAACAACACTCGTCGGCAGCGTC"
$ws.Range("B93").Value = "F0572-TTGTCTCCTG"
$ws.Range("C93").Value = "AATGATACGGCGACCACCGAGATCTACACTTGTCTCCTGTCGTCGGCAGCGTC"
$ws.Range("B94").Value = "F0573-TACAACTTCG"
$ws.Range("C94").Value = "AATGATACGGCGACCACCGAGATCTACACTACAACTTCGTCGTCGGCAGCGTC"
$ws.Range("B95").Value = "F0574-TTCAGTGAGG"
$ws.Range("C95").Value = "AATGATACGGCGACCACCGAGATCTACACTTCAGTGAGGTCGTCGGCAGCGTC"
$ws.Range("B96").Value = "F0575-CTTGGATCCT"
$ws.Range("C96").Value = "AATGATACGGCGACCACCGAGATCTACACCTTGGATCCTTCGTCGGCAGCGTC"
$ws.Range("B97").Value = "F0576-CATGCTACGA"
$ws.Range("C97").Value = "AATGATACGGCGACCACCGAGATCTACACCATGCTACGATCGTCGGCAGCGTC"
